$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "enrollment device" records (rows 6-7, id 589 and 638 —
# Vostro/Dell/DKS "To take enrollments" in English + Arabic). Deleting
# these two rows shifts the remaining device rows up so row 8 (id 736)
# becomes row 6, etc.
$ws.Rows("6:7").Delete()

# Match the page setup recorded for the sheet (A4-ish defaults plus an
# explicit paper size/orientation, as captured by the print settings).
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$ws.Range("E16").Select()
